# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# --- Rename the "Include from SetOperator" sheet to "Include #0" ---
$include = $wb.Worksheets.Item("Include from SetOperator")
$include.Name = "Include #0"

# --- Update the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Bump version + date
$meta.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10),
# copying the formatting of the existing data rows so the inserted
# row matches the rest of the table's style.
$meta.Rows.Item(11).Insert()
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Cells.Item(11, 1).Value = "Jurisdiction"
$meta.Cells.Item(11, 2).Value = ""
